# Update "想去人数" (F column) counts and, where a listing became sellable
# again, the "最低票价" (G column) value on the 北京-漫展信息 workbook.
#
# Sheet order in the workbook:
#   1 = 展览      (Exhibitions)
#   2 = 演出      (Performances)
#   3 = 本地生活  (Local life)
#   4 = 全部类型  (All types - aggregates the above three sheets)

$wb = $excel.ActiveWorkbook

# ---- 展览 (sheet 1) ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(6, 6).Value = 277
$ws.Cells.Item(7, 6).Value = 9299
$ws.Cells.Item(8, 6).Value = 234
$ws.Cells.Item(9, 6).Value = 25
$ws.Cells.Item(10, 6).Value = 673
$ws.Cells.Item(11, 6).Value = 1799
$ws.Cells.Item(12, 6).Value = 37
$ws.Cells.Item(13, 6).Value = 88
$ws.Cells.Item(14, 6).Value = 2466
$ws.Cells.Item(15, 6).Value = 116
$ws.Cells.Item(16, 6).Value = 3752
$ws.Cells.Item(17, 6).Value = 272
$ws.Cells.Item(18, 6).Value = 108
$ws.Cells.Item(19, 6).Value = 118
$ws.Cells.Item(21, 6).Value = 229
$ws.Cells.Item(22, 6).Value = 186
$ws.Cells.Item(24, 6).Value = 51
$ws.Cells.Item(25, 6).Value = 246
$ws.Cells.Item(26, 6).Value = 497
$ws.Cells.Item(27, 6).Value = 107
$ws.Cells.Item(28, 6).Value = 1076
$ws.Cells.Item(29, 6).Value = 435
$ws.Cells.Item(30, 6).Value = 4294
$ws.Cells.Item(31, 6).Value = 65
$ws.Cells.Item(32, 6).Value = 44
$ws.Cells.Item(32, 7).Value = "'70"
$ws.Cells.Item(33, 6).Value = 257
$ws.Cells.Item(34, 6).Value = 32
$ws.Cells.Item(34, 7).Value = "'70"

# ---- 演出 (sheet 2) ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 40
$ws.Cells.Item(3, 6).Value = 34
$ws.Cells.Item(5, 6).Value = 14

# ---- 本地生活 (sheet 3) ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 6).Value = 939

# ---- 全部类型 (sheet 4) ----
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 6).Value = 939
$ws.Cells.Item(5, 6).Value = 40
$ws.Cells.Item(6, 6).Value = 34
$ws.Cells.Item(10, 6).Value = 277
$ws.Cells.Item(11, 6).Value = 9299
$ws.Cells.Item(12, 6).Value = 234
$ws.Cells.Item(13, 6).Value = 25
$ws.Cells.Item(14, 6).Value = 673
$ws.Cells.Item(15, 6).Value = 1799
$ws.Cells.Item(16, 6).Value = 37
$ws.Cells.Item(17, 6).Value = 88
$ws.Cells.Item(19, 6).Value = 2466
$ws.Cells.Item(20, 6).Value = 116
$ws.Cells.Item(21, 6).Value = 3752
$ws.Cells.Item(22, 6).Value = 272
$ws.Cells.Item(23, 6).Value = 108
$ws.Cells.Item(24, 6).Value = 118
$ws.Cells.Item(26, 6).Value = 229
$ws.Cells.Item(27, 6).Value = 186
$ws.Cells.Item(28, 6).Value = 14
$ws.Cells.Item(30, 6).Value = 51
$ws.Cells.Item(31, 6).Value = 246
$ws.Cells.Item(32, 6).Value = 497
$ws.Cells.Item(33, 6).Value = 107
$ws.Cells.Item(34, 6).Value = 1076
$ws.Cells.Item(35, 6).Value = 435
$ws.Cells.Item(36, 6).Value = 4294
$ws.Cells.Item(37, 6).Value = 65
$ws.Cells.Item(38, 6).Value = 44
$ws.Cells.Item(38, 7).Value = "'70"
$ws.Cells.Item(39, 6).Value = 257
$ws.Cells.Item(40, 6).Value = 32
$ws.Cells.Item(40, 7).Value = "'70"
